$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; this shifts the existing rows 14-27
# down to 15-28, carrying their data/formatting with them.
$ws.Rows("14:14").Insert()

# Fill the newly inserted row 14 with the new weekly record.
$ws.Cells.Item(14, 1).Value()  = 1
$ws.Cells.Item(14, 2).Value()  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14, 3).Value()  = "Arica y Parinacota"
$ws.Cells.Item(14, 4).Value()  = 44601
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14, 5).Value()  = 15
$ws.Cells.Item(14, 6).Value()  = 100112044
$ws.Cells.Item(14, 7).Value()  = "Perejil"
$ws.Cells.Item(14, 8).Value()  = "Sin especificar"
$ws.Cells.Item(14, 9).Value()  = "Primera"
$ws.Cells.Item(14, 10).Value() = 270
$ws.Cells.Item(14, 11).Value() = 2200
$ws.Cells.Item(14, 12).Value() = 2500
$ws.Cells.Item(14, 13).Value() = 2350
$ws.Cells.Item(14, 14).Value() = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(14, 15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value() = 1175
$ws.Cells.Item(14, 17).Value() = 2
$ws.Cells.Item(14, 18).Value() = "Hortaliza"
